# Workbook "Uren log.xlsx" update:
#  - add a new logged-hours entry on the "Thomas" sheet (row 21):
#      "Finishing Maersk scraper", 18-11-2022, 2 hours
#  - extend the running-total SUM formula in E3 to cover the new rows
#  - leave the active selection on H11 (matches the saved view state)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Thomas")

# Copy the date formatting from the row above so the new date cell keeps
# the existing "short date" style instead of Excel inventing a new one.
$ws.Range("B20").Copy()
$ws.Range("B21").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false

# New row of data
$ws.Range("A21").Value = "Finishing Maersk scraper"
$ws.Range("B21").Value = [DateTime]"2022-11-18"
$ws.Range("C21").Value = 2

# Extend the totals formula to include the newly added rows
$ws.Range("E3").Formula = "=SUM(C2:C30)"

# Restore the selection that was active when the file was saved
$ws.Range("H11").Select()
